# Edit sheet Card24 by admin
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Row 2: Event / Correction replaced, Serviced by filled in
$ws.Range("M2").Value = "تم استبدال بواحد جديد"
$ws.Range("N2").Value = "تم استبدال بواحد جديد"
$ws.Range("O2").Value = "م.صيام"

# Row 3: "Servised by" cell gets the literal placeholder text "nan"
$ws.Range("O3").Value = "nan"

# Row 4: Event note added, Correction cleared out
$ws.Range("M4").Value = "ليكر ان مكسور"
$ws.Range("N4").Value = ""

# Rows 5-12: "Servised by" cells get the literal placeholder text "nan"
$ws.Range("O5").Value = "nan"
$ws.Range("O6").Value = "nan"
$ws.Range("O7").Value = "nan"
$ws.Range("O8").Value = "nan"
$ws.Range("O9").Value = "nan"
$ws.Range("O10").Value = "nan"
$ws.Range("O11").Value = "nan"
$ws.Range("O12").Value = "nan"
